$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.681.50"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "1.757.05"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.93"
$ws.Range("E5").Value = "  +0.61%  "
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4591"
$ws.Range("E7").Value = "  +7.66%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3595"
$ws.Range("E8").Value = "  -1.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07521"
$ws.Range("E9").Value = "  +0.95%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.13"
$ws.Range("E10").Value = "  -2.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.102"
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.77"
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("E14").Value = "  -1.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.118"
$ws.Range("E15").Value = "  -3.24%  "
$ws.Range("D16").Value = "1.758.43"
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.47"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001068"
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06412"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("E21").Value = "  -1.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.829"
$ws.Range("E22").Value = "  -2.34%  "
$ws.Range("D23").Value = "27.742.50"
$ws.Range("E23").Value = "  -0.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.27"
$ws.Range("E24").Value = "  -0.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.100"
$ws.Range("E25").Value = "  +1.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.19"
$ws.Range("E26").Value = "  +3.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.36"
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("D28").Value = "1.960.68"
$ws.Range("E28").Value = "  -1.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.085"
$ws.Range("E29").Value = "  -3.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.49"
$ws.Range("E30").Value = "  +0.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.060"
$ws.Range("E31").Value = "  -7.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09226"
$ws.Range("E32").Value = "  +3.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.671"
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.558"
$ws.Range("E34").Value = "  -0.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.94"
$ws.Range("E35").Value = "  -3.74%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02310"
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2107"
$ws.Range("E37").Value = "  -0.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06060"
$ws.Range("E38").Value = "  +0.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6356"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.987"
$ws.Range("E40").Value = "  -0.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.209"
$ws.Range("E41").Value = "  +2.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.378"
$ws.Range("E42").Value = "  -1.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.809"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.31"
$ws.Range("E44").Value = "  -0.88%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5920"
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("E46").Value = "  +0.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "123.61"
$ws.Range("E47").Value = "  +0.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.951"
$ws.Range("E48").Value = "  -2.79%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.149"
$ws.Range("E49").Value = "  -2.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06868"
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.26"
$ws.Range("E51").Value = "  -2.65%  "
